# "Klar til nattkjøring fortsettelse" - append the next batch of simulation
# results (rows 12-14) produced by the overnight run to the "Simulation"
# sheet, right below the last existing result row (11), and leave the
# selection where the analyst left off (T20).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New result row 12 ------------------------------------------------
$ws.Cells.Item(12, 1).Value = 0
$ws.Cells.Item(12, 2).Value = 0
$ws.Cells.Item(12, 3).Value = -0.89999999999999991
$ws.Cells.Item(12, 4).Value = 0.10000000000000009
$ws.Cells.Item(12, 5).Value = 0.6
$ws.Cells.Item(12, 6).Value = 0.3
$ws.Cells.Item(12, 7).Value = 0.1
$ws.Cells.Item(12, 8).Value = 0.6
$ws.Cells.Item(12, 9).Value = 0.4
$ws.Cells.Item(12, 10).Value = 5
$ws.Cells.Item(12, 11).Value = 18
$ws.Cells.Item(12, 12).Value = "HEURISTIC_VERSION_2"
$ws.Cells.Item(12, 13).Value = "EVERY_VEHICLE_ARRIVAL"
$ws.Cells.Item(12, 14).Value = 1
$ws.Cells.Item(12, 15).Value = 20
$ws.Cells.Item(12, 16).Value = 7
$ws.Cells.Item(12, 17).Value = 11
$ws.Cells.Item(12, 18).Value = 5
$ws.Cells.Item(12, 19).Value = 1
$ws.Cells.Item(12, 20).Value = 3
$ws.Cells.Item(12, 21).Value = 3
$ws.Cells.Item(12, 22).Value = 15
$ws.Cells.Item(12, 23).Value = 5
$ws.Cells.Item(12, 24).Value = 3.8678629146118579
$ws.Cells.Item(12, 25).Value = 3.8918107891778635
$ws.Cells.Item(12, 26).Value = 50.231950857711084
$ws.Cells.Item(12, 27).Value = 28.733333333333334
$ws.Cells.Item(12, 28).Value = 8.5134273967478027
$ws.Range("A12:AB12").Style = "Normal"

# --- New result row 13 ------------------------------------------------
$ws.Cells.Item(13, 1).Value = 0
$ws.Cells.Item(13, 2).Value = 0
$ws.Cells.Item(13, 3).Value = -0.99999999999999989
$ws.Cells.Item(13, 4).Value = 0
$ws.Cells.Item(13, 5).Value = 0.6
$ws.Cells.Item(13, 6).Value = 0.3
$ws.Cells.Item(13, 7).Value = 0.1
$ws.Cells.Item(13, 8).Value = 0.6
$ws.Cells.Item(13, 9).Value = 0.4
$ws.Cells.Item(13, 10).Value = 5
$ws.Cells.Item(13, 11).Value = 18
$ws.Cells.Item(13, 12).Value = "HEURISTIC_VERSION_2"
$ws.Cells.Item(13, 13).Value = "EVERY_VEHICLE_ARRIVAL"
$ws.Cells.Item(13, 14).Value = 1
$ws.Cells.Item(13, 15).Value = 20
$ws.Cells.Item(13, 16).Value = 7
$ws.Cells.Item(13, 17).Value = 11
$ws.Cells.Item(13, 18).Value = 5
$ws.Cells.Item(13, 19).Value = 1
$ws.Cells.Item(13, 20).Value = 3
$ws.Cells.Item(13, 21).Value = 3
$ws.Cells.Item(13, 22).Value = 15
$ws.Cells.Item(13, 23).Value = 5
$ws.Cells.Item(13, 24).Value = 4.3783958109319006
$ws.Cells.Item(13, 25).Value = 4.401511680514826
$ws.Cells.Item(13, 26).Value = 50.242622762003911
$ws.Cells.Item(13, 27).Value = 31.333333333333332
$ws.Cells.Item(13, 28).Value = 7.7750928836836097
$ws.Range("A13:AB13").Style = "Normal"

# --- New result row 14 ------------------------------------------------
$ws.Cells.Item(14, 1).Value = 0
$ws.Cells.Item(14, 2).Value = 0.1
$ws.Cells.Item(14, 3).Value = 0
$ws.Cells.Item(14, 4).Value = 0.9
$ws.Cells.Item(14, 5).Value = 0.6
$ws.Cells.Item(14, 6).Value = 0.3
$ws.Cells.Item(14, 7).Value = 0.1
$ws.Cells.Item(14, 8).Value = 0.6
$ws.Cells.Item(14, 9).Value = 0.4
$ws.Cells.Item(14, 10).Value = 5
$ws.Cells.Item(14, 11).Value = 18
$ws.Cells.Item(14, 12).Value = "HEURISTIC_VERSION_2"
$ws.Cells.Item(14, 13).Value = "EVERY_VEHICLE_ARRIVAL"
$ws.Cells.Item(14, 14).Value = 1
$ws.Cells.Item(14, 15).Value = 20
$ws.Cells.Item(14, 16).Value = 7
$ws.Cells.Item(14, 17).Value = 11
$ws.Cells.Item(14, 18).Value = 5
$ws.Cells.Item(14, 19).Value = 1
$ws.Cells.Item(14, 20).Value = 3
$ws.Cells.Item(14, 21).Value = 3
$ws.Cells.Item(14, 22).Value = 15
$ws.Cells.Item(14, 23).Value = 5
$ws.Cells.Item(14, 24).Value = 0.7990273316982297
$ws.Cells.Item(14, 25).Value = 0.80891343188631126
$ws.Cells.Item(14, 26).Value = 50.275755849011354
$ws.Cells.Item(14, 27).Value = 20.066666666666666
$ws.Cells.Item(14, 28).Value = 12.343588938923956
$ws.Range("A14:AB14").Style = "Normal"

# Leave the cursor parked on T20, ready for the next night's run.
$ws.Range("T20").Select()
